# "Add files via upload" — appends one new product row (row 39) to Sheet1
# of HARGA ONLINE.xlsx: "TASBIH KAYU BUTIR", priced Rp4,750, with "-"
# placeholders in the paper-type/size columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- values for the new row ------------------------------------------
$ws.Range("A39").Value = 38
$ws.Range("B39").Value = "TASBIH KAYU BUTIR"
$ws.Range("C39").Value = "-"
$ws.Range("D39").Value = "-"
$ws.Range("E39").Value = 4750

$dataRange = $ws.Range("A39:D39")

# ---- font: match the body rows above (Calibri 12, not bold) ----------
$dataRange.Font.Size = 12

# ---- borders: thin left/right on A39:D39 (no top/bottom) -------------
foreach ($col in @("A", "B", "C", "D")) {
    $cell = $ws.Range($col + "39")
    $cell.Borders(7).LineStyle = 1
    $cell.Borders(7).Weight = 2
    $cell.Borders(10).LineStyle = 1
    $cell.Borders(10).Weight = 2
}

# ---- alignment ---------------------------------------------------------
# -4108 = xlCenter, -4131 = xlLeft
$ws.Range("A39").HorizontalAlignment = -4108
$ws.Range("A39").VerticalAlignment = -4108
$ws.Range("B39").HorizontalAlignment = -4131
$ws.Range("C39").HorizontalAlignment = -4108
$ws.Range("D39").HorizontalAlignment = -4108
$ws.Range("E39").HorizontalAlignment = -4108

# ---- number format for the price cell ---------------------------------
$ws.Range("E39").NumberFormat = """Rp""#,##0;[Red]\-""Rp""#,##0"

# ---- row height, matching the other data rows (15.75pt) ---------------
$ws.Rows("39").RowHeight = 15.75

# ---- update the view: scroll down, select C44 as the active cell -------
$ws.Range("C44").Select()
